# Allow all plant types to be built for reliability and set dispatch cost
# multiplier to 1 (#232)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPTBfRN")

# Set the Boolean "built for reliability" flags to 1 for rows 2-11
# (hard coal, natural gas steam turbine, natural gas combined cycle,
# nuclear, hydro, onshore wind, solar PV, solar thermal, biomass, geothermal)
$ws.Range("B2:B11").Value = 1
